$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.917.03"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.777.41"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'598.72"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'163.23"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'6.58"
$ws.Range("E11").Value = "  +4.34%  "
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "'35.34"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").Value = "4.409.68"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.775.21"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "67.887.69"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'18.26"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "'7.00"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "'457.80"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "'9.58"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'82.67"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'0.0000142"
$ws.Range("E24").Value = "  -6.41%  "
$ws.Range("D25").Value = "'11.91"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "3.927.76"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -7.98%  "
$ws.Range("D33").Value = "'28.97"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'8.91"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "'3.17"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'43.43"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "'47.25"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "'152.64"
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").Value = "'8.28"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "'386.48"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").Value = "'26.39"
$ws.Range("E51").Value = "  -3.67%  "
